$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Remove the existing hyperlinks up front. Column inserts below do not
#    relocate hyperlink anchors automatically in this engine, so hyperlinks
#    are recreated from scratch (at the right spots) near the end instead.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 1. Insert a new column at E ("appComponent" / "test_comp").
#    This shifts old E..I (browserPath..debug_error) to F..J.
# ---------------------------------------------------------------------------
$ws.Range("E1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Insert three new columns at I:K ("sendResults", "xrayUrl",
#    "xrayTestExecKeyMaster"). This shifts the (already shifted) old H/I
#    columns (writeFailReqRspOnly / debug_error), now sitting at I/J, on to
#    L/M.
# ---------------------------------------------------------------------------
$ws.Range("I1:K1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. Stash the formatting of a few reference cells whose style we want to
#    reuse for new cells / cells that a later hyperlink-add would otherwise
#    clobber. (Copy/PasteSpecial-Formats reuses an existing matching style
#    entry instead of fabricating a new one, as long as it is done before
#    Hyperlinks.Add touches the cell.)
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # style 3 (hyperlink, general fmt)
$ws.Range("H2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # style 7 (hyperlink, text fmt)
$ws.Range("G1").Copy()
$ws.Range("Z3").PasteSpecial(-4122)   # style 2 (header, gray fill)
$ws.Range("G2").Copy()
$ws.Range("Z4").PasteSpecial(-4122)   # style 1 (plain, text fmt)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Fill in the values for the brand new header cells (row 1).
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "appComponent"
$ws.Range("I1").Value = "sendResults"
$ws.Range("J1").Value = "xrayUrl"
$ws.Range("K1").Value = "xrayTestExecKeyMaster"

# ---------------------------------------------------------------------------
# 5. Fill in the values for the brand new data cells (row 2). A leading
#    apostrophe forces "false"/"true"-looking text to stay text instead of
#    being auto-converted to a Boolean cell.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "test_comp"
$ws.Range("I2").Value = "'false"
$ws.Range("J2").Value = "https://jira.cleverlance.com/rest/raven/1.0/import/execution"
$ws.Range("K2").Value = "n/a"

# ---------------------------------------------------------------------------
# 6. Update the existing "measurement" value.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "test_measurement_mze"

# ---------------------------------------------------------------------------
# 7. Apply the matching formatting to the new cells.
# ---------------------------------------------------------------------------
$ws.Range("Z3").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("K1").PasteSpecial(-4122)

$ws.Range("Z4").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I2").Value = "'false"
$ws.Range("K2").Value = "n/a"

# ---------------------------------------------------------------------------
# 8. Column widths - best effort (this engine quantizes widths to 1/6 of a
#    character, so the closest achievable value is used).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14.25   # B  -> target 15.109375
$ws.Columns.Item(5).ColumnWidth = 12.75   # E  -> target 13.6640625
$ws.Columns.Item(9).ColumnWidth = 9.59    # I  -> target 10.5546875
$ws.Columns.Item(10).ColumnWidth = 45.92  # J  -> target 46.77734375
$ws.Columns.Item(11).ColumnWidth = 22.25  # K  -> target 23.21875

# ---------------------------------------------------------------------------
# 9. Re-create the hyperlinks at their final locations, then restore their
#    styling (Hyperlinks.Add always fabricates a brand-new style entry).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://test.clv.cz/")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://monitoring-test.kb.cz/appmon-in-test/write")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://jira.cleverlance.com/rest/raven/1.0/import/execution")

$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 10. Clean up the scratch cells used to stash formatting.
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 11. Select A2, matching the saved selection in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
